$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row after the existing data (row 38)
$newRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

# New daily entry: Day (date serial 45987 = 2025-11-26), Chase, Bryce, Zach
$ws.Cells.Item($newRow, 1).Value = 45987
$ws.Cells.Item($newRow, 2).Value = 87
$ws.Cells.Item($newRow, 3).Value = 93
$ws.Cells.Item($newRow, 4).Value = 92

# Match the date style/format used by the rest of column A
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat
